$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: new Debit formula
$ws.Range("D3").Formula = "=45000+405000"

# Row 4: TRANSFER BCA
$ws.Range("B4").Value = "TRANSFER BCA"
$ws.Range("D4").Formula = "=8510000+12950000+540000+1385000+975000+800000+140000+6118000"

# Row 5: TAX - P.Tata
$ws.Range("B5").Value = "TAX - P.Tata"
$ws.Range("D5").Formula = "=200000"

# Row 6: TAX - Iuran ARIESTA
$ws.Range("B6").Value = "TAX - Iuran ARIESTA"
$ws.Range("D6").Formula = "=660000"

# Row 7: prive
$ws.Range("B7").Value = "prive"
$ws.Range("D7").Value = 2000000

# Row 8: SALES - cash/retail
$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("C8").Formula = "=27600475+23611525+8510000-50104000"

# Row 9: A/R
$ws.Range("B9").Value = "A/R"
$ws.Range("C9").Formula = "=3280000+2838000+27571000+37195000+50104000"

# Row 10: TRANSFER BCA AA
$ws.Range("B10").Value = "TRANSFER BCA AA"
$ws.Range("D10").Formula = "=64855000"

# Row 11: SELISIH - lebih
$ws.Range("B11").Value = "SELISIH - lebih"
$ws.Range("C11").Value = 80000

# Row 12: SETOR KE BANK
$ws.Range("B12").Value = "SETOR KE BANK"
$ws.Range("D12").Value = 31000000

# Row 13: SOLAR - kijang
$ws.Range("B13").Value = "SOLAR - kijang"
$ws.Range("D13").Value = 300000

# Row 14: new date entry
$ws.Range("A14").Value = 44208

# Update view: frozen pane top-left cell and active selection
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B34").Select()
